{"js": "// Remove the two To-Do bullet items that were dropped from the list:\n//  - \"Minimum grain size to insert twin user input (or rule)\"\n//  - \"Check twin thickness per grain to make sure its thicker than 1 voxel\"\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst targets = [\n  \"Minimum grain size to insert twin user input (or rule)\",\n  \"Check twin thickness per grain to make sure its thicker than 1 voxel\"\n];\n\nfor (let i = paragraphs.items.length - 1; i >= 0; i--) {\n  const text = paragraphs.items[i].text.trim();\n  if (targets.indexOf(text) !== -1) {\n    paragraphs.items[i].delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the two To-Do bullet items that were dropped from the list:\n#  - \"Minimum grain size to insert twin user input (or rule)\"\n#  - \"Check twin thickness per grain to make sure its thicker than 1 voxel\"\n$d = $word.ActiveDocument\n\n$targets = @(\n  \"Minimum grain size to insert twin user input (or rule)\",\n  \"Check twin thickness per grain to make sure its thicker than 1 voxel\"\n)\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.Trim()\n    if ($targets -contains $t) {\n        $p.Range.Delete()\n    }\n}\n"}
